$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlineShape($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footer (default, rId8 -> footer2.xml): Pearson logo image2.png -> image1.png
$defaultFooter = $sec.Footers.Item(1)
Rename-InlineShape $defaultFooter.Range.InlineShapes.Item(1) "image1.png"

# Footer (first page, rId9 -> footer1.xml): Pearson logo image2.png -> image1.png
$firstFooter = $sec.Footers.Item(2)
Rename-InlineShape $firstFooter.Range.InlineShapes.Item(1) "image1.png"

# Header (first page, rId7 -> header1.xml): BTEC logo image1.jpg -> image2.jpg
$firstHeader = $sec.Headers.Item(2)
Rename-InlineShape $firstHeader.Range.InlineShapes.Item(1) "image2.jpg"
